# Update field names in the Statement_of_Receivable-Customer template.
# Rows 6-8 had the {ACRefCode...} / {PartnerCode...} merge-field placeholders
# swapped between columns B and C; this corrects them so column B carries
# the PartnerCode family and column C carries the ACRefCode family.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value = "{PartnerCode}"
$ws.Range("C6").Value = "{ACRefCode}"

$ws.Range("B7").Value = "{PartnerCodeDt}"
$ws.Range("C7").Value = "{ACRefCodeDt}"

$ws.Range("B8").Value = "{PartnerCodeAdv}"
$ws.Range("C8").Value = "{ACRefCodeAdv}"
